{"js": "// Rename the hill-parameter labels in the algorithm list to lowercase,\n// snake_case identifiers:\n//   \"Height = 46\"            -> \"height = 46\"\n//   \"Points per meter = 2\"   -> \"points_per_meter = 2\"\n//   \"Par (distance) = 90\"    -> \"par (distance) = 90\"\n//   \"Height = 70\"            -> \"height = 70\"\n//   \"Points per meter = 1.8\" -> \"points_per_meter = 1.8\"\n//   \"Par (distance) = 120\"   -> \"par (distance) = 120\"\n\nconst replacements = [\n  [\"Height = 46\", \"height = 46\"],\n  [\"Points per meter = 2\", \"points_per_meter = 2\"],\n  [\"Par (distance) = 90\", \"par (distance) = 90\"],\n  [\"Height = 70\", \"height = 70\"],\n  [\"Points per meter = 1.8\", \"points_per_meter = 1.8\"],\n  [\"Par (distance) = 120\", \"par (distance) = 120\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the hill-parameter labels in the algorithm list to lowercase,\n# snake_case identifiers:\n#   \"Height = 46\"            -> \"height = 46\"\n#   \"Points per meter = 2\"   -> \"points_per_meter = 2\"\n#   \"Par (distance) = 90\"    -> \"par (distance) = 90\"\n#   \"Height = 70\"            -> \"height = 70\"\n#   \"Points per meter = 1.8\" -> \"points_per_meter = 1.8\"\n#   \"Par (distance) = 120\"   -> \"par (distance) = 120\"\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Height = 46\", \"height = 46\"),\n    @(\"Points per meter = 2\", \"points_per_meter = 2\"),\n    @(\"Par (distance) = 90\", \"par (distance) = 90\"),\n    @(\"Height = 70\", \"height = 70\"),\n    @(\"Points per meter = 1.8\", \"points_per_meter = 1.8\"),\n    @(\"Par (distance) = 120\", \"par (distance) = 120\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
